# Auto-generated data correction pass (scheduled Sheets -> Kujata_Profits sync).
# Re-applies the refreshed currentAveragePrice*/LevePrice*/LeveProfit* figures
# for each Leve table (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the latest pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 2730.0747
$ws.Range("I15").Value = 2730.0747
$ws.Range("K15").Value = 8190.2241
$ws.Range("M15").Value = -8021.2241
# Row 19
$ws.Range("H19").Value = 899.5
$ws.Range("I19").Value = 799
$ws.Range("J19").Value = 933
$ws.Range("K19").Value = 799
$ws.Range("L19").Value = 933
$ws.Range("M19").Value = -624
$ws.Range("N19").Value = -1283
# Row 100
$ws.Range("H100").Value = 2267.4827
$ws.Range("I100").Value = 2250.2593
$ws.Range("J100").Value = 2500
$ws.Range("K100").Value = 2250.2593
$ws.Range("L100").Value = 2500
$ws.Range("M100").Value = -1709.2593
$ws.Range("N100").Value = -3582
# Row 107
$ws.Range("H107").Value = 2873.1765
$ws.Range("I107").Value = 3822.5
$ws.Range("J107").Value = 2355.3635
$ws.Range("K107").Value = 3822.5
$ws.Range("L107").Value = 2355.3635
$ws.Range("M107").Value = -1902.5
$ws.Range("N107").Value = -6195.363499999999
# Row 112
$ws.Range("H112").Value = 1879.4348
$ws.Range("J112").Value = 2052.6667
$ws.Range("L112").Value = 6158.000100000001
$ws.Range("N112").Value = -8374.000100000001
# Row 129
$ws.Range("H129").Value = 648.53656
$ws.Range("I129").Value = 333.125
$ws.Range("J129").Value = 850.4
$ws.Range("K129").Value = 999.375
$ws.Range("L129").Value = 2551.2
$ws.Range("M129").Value = 4000.625
$ws.Range("N129").Value = -12551.2
# Row 137
$ws.Range("H137").Value = 1198.4717
$ws.Range("I137").Value = 803.56665
$ws.Range("J137").Value = 1713.5652
$ws.Range("K137").Value = 2410.69995
$ws.Range("L137").Value = 5140.6956
$ws.Range("M137").Value = 139.3000499999998
$ws.Range("N137").Value = -10240.6956
# Row 138
$ws.Range("H138").Value = 1314.0465
$ws.Range("J138").Value = 1654.2916
$ws.Range("L138").Value = 4962.8748
$ws.Range("N138").Value = -15242.8748
# Row 141
$ws.Range("H141").Value = 587.3461
$ws.Range("I141").Value = 530.84
$ws.Range("K141").Value = 1592.52
$ws.Range("M141").Value = 3587.48

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3915.074
$ws.Range("I32").Value = 3555.1216
$ws.Range("J32").Value = 7720.2856
$ws.Range("K32").Value = 3555.1216
$ws.Range("L32").Value = 7720.2856
$ws.Range("M32").Value = -3268.1216
$ws.Range("N32").Value = -8294.285599999999
# Row 74
$ws.Range("H74").Value = 1097.6666
$ws.Range("I74").Value = 770.2963
$ws.Range("K74").Value = 770.2963
$ws.Range("M74").Value = 103.7037
# Row 77
$ws.Range("H77").Value = 1097.6666
$ws.Range("I77").Value = 770.2963
$ws.Range("K77").Value = 3851.4815
$ws.Range("M77").Value = 516.5185000000001
# Row 132
$ws.Range("H132").Value = 1479.3611
$ws.Range("I132").Value = 1548.091
$ws.Range("J132").Value = 1371.3572
$ws.Range("K132").Value = 4644.272999999999
$ws.Range("L132").Value = 4114.071599999999
$ws.Range("M132").Value = -2114.272999999999
$ws.Range("N132").Value = -9174.071599999999

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 1177.05
$ws.Range("I107").Value = 1041.2222
$ws.Range("J107").Value = 2399.5
$ws.Range("K107").Value = 1041.2222
$ws.Range("L107").Value = 2399.5
$ws.Range("M107").Value = 878.7778000000001
$ws.Range("N107").Value = -6239.5
# Row 132
$ws.Range("H132").Value = 59499.75
$ws.Range("J132").Value = 59499.75
$ws.Range("L132").Value = 59499.75
$ws.Range("N132").Value = -69619.75
# Row 134
$ws.Range("H134").Value = 5734.885
$ws.Range("I134").Value = 1309
$ws.Range("K134").Value = 3927
$ws.Range("M134").Value = -1392

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1720.3125
$ws.Range("I31").Value = 1549.4524
$ws.Range("K31").Value = 1549.4524
$ws.Range("M31").Value = -1254.4524
# Row 34
$ws.Range("H34").Value = 1720.3125
$ws.Range("I34").Value = 1549.4524
$ws.Range("K34").Value = 1549.4524
$ws.Range("M34").Value = -1347.4524
# Row 58
$ws.Range("H58").Value = 823.45
$ws.Range("I58").Value = 733.6111
$ws.Range("K58").Value = 733.6111
$ws.Range("M58").Value = -530.6111
# Row 132
$ws.Range("H132").Value = 4641.3887
$ws.Range("I132").Value = 5360.724
$ws.Range("J132").Value = 1661.2858
$ws.Range("K132").Value = 16082.172
$ws.Range("L132").Value = 4983.857400000001
$ws.Range("M132").Value = -13552.172
$ws.Range("N132").Value = -10043.8574
# Row 133
$ws.Range("H133").Value = 36596.5
$ws.Range("J133").Value = 36596.5
$ws.Range("L133").Value = 36596.5
$ws.Range("N133").Value = -41656.5
# Row 136
$ws.Range("H136").Value = 823.45
$ws.Range("I136").Value = 733.6111
$ws.Range("K136").Value = 2200.8333
$ws.Range("M136").Value = 349.1667000000002

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 615.3333
$ws.Range("J113").Value = 654.0833
$ws.Range("L113").Value = 1962.2499
$ws.Range("N113").Value = -6302.2499
# Row 122
$ws.Range("H122").Value = 778.9355
$ws.Range("I122").Value = 674.46155
$ws.Range("J122").Value = 854.3889
$ws.Range("K122").Value = 6070.15395
$ws.Range("L122").Value = 7689.5001
$ws.Range("M122").Value = -3620.15395
$ws.Range("N122").Value = -12589.5001
# Row 131
$ws.Range("H131").Value = 31251500
$ws.Range("J131").Value = 1861.174
$ws.Range("L131").Value = 5583.522
$ws.Range("N131").Value = -15663.522
# Row 132
$ws.Range("H132").Value = 1356.6364
$ws.Range("I132").Value = 987.9231
$ws.Range("J132").Value = 1889.2222
$ws.Range("K132").Value = 8891.3079
$ws.Range("L132").Value = 17002.9998
$ws.Range("M132").Value = -6361.3079
$ws.Range("N132").Value = -22062.9998
# Row 140
$ws.Range("H140").Value = 26629.38
$ws.Range("J140").Value = 3037.2173
$ws.Range("L140").Value = 9111.651899999999
$ws.Range("N140").Value = -19471.6519

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 2368.5557
$ws.Range("I132").Value = 2056.4614
$ws.Range("J132").Value = 3180
$ws.Range("K132").Value = 6169.3842
$ws.Range("L132").Value = 9540
$ws.Range("M132").Value = -3639.3842
$ws.Range("N132").Value = -14600

$ws = $wb.Worksheets.Item("LTW")
# Row 80
$ws.Range("H80").Value = 16000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 16000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 16000
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -18246
# Row 83
$ws.Range("H83").Value = 16000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 16000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 48000
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -59232
# Row 131
$ws.Range("H131").Value = 42500
$ws.Range("J131").Value = 42500
$ws.Range("L131").Value = 42500
$ws.Range("N131").Value = -52580
# Row 132
$ws.Range("H132").Value = 42160.24
$ws.Range("I132").Value = 2172.9092
$ws.Range("J132").Value = 73578.86
$ws.Range("K132").Value = 6518.7276
$ws.Range("L132").Value = 220736.58
$ws.Range("M132").Value = -3988.7276
$ws.Range("N132").Value = -225796.58

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 381.83334
$ws.Range("I81").Value = 381.83334
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 763.66668
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = 297.33332
$ws.Range("N81").ClearContents()
# Row 84
$ws.Range("H84").Value = 381.83334
$ws.Range("I84").Value = 381.83334
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 3818.3334
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 1485.6666
$ws.Range("N84").ClearContents()
# Row 124
$ws.Range("H124").Value = 64500
$ws.Range("J124").Value = 64500
$ws.Range("L124").Value = 64500
$ws.Range("N124").Value = -74320
# Row 131
$ws.Range("H131").Value = 49250
$ws.Range("J131").Value = 49250
$ws.Range("L131").Value = 49250
$ws.Range("N131").Value = -59330
# Row 136
$ws.Range("H136").Value = 580.4761999999999
$ws.Range("I136").Value = 558.75
$ws.Range("J136").Value = 650
$ws.Range("K136").Value = 1676.25
$ws.Range("L136").Value = 1950
$ws.Range("M136").Value = 873.75
$ws.Range("N136").Value = -7050
